# Add 7 new rows (422-428) of landscaping data collected on 7/9/2025
# (date serial 45847), matching the rows appended in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Write the plain (non-formula) cell values for the new rows.
# ---------------------------------------------------------------------

# Row 422
$ws.Range("A422").Value = 45847
$ws.Range("B422").Value = "Flowering"
$ws.Range("C422").Value = "Large"
$ws.Range("D422").Value = 66
$ws.Range("E422").Value = 77
$ws.Range("G422").Value = 0.25
$ws.Range("H422").Value = 0.1
$ws.Range("I422").Value = "No"
$ws.Range("J422").Value = 2
$ws.Range("K422").Value = "Dark"
$ws.Range("L422").Value = 5
$ws.Range("M422").Value = 0.84
$ws.Range("N422").Value = 70
$ws.Range("O422").Value = 30.01
$ws.Range("P422").Value = 7
$ws.Range("Q422").Value = 0.86
$ws.Range("R422").Value = 5.8
$ws.Range("S422").Value = 57
$ws.Range("T422").Value = 0

# Row 423
$ws.Range("A423").Value = 45847
$ws.Range("B423").Value = "Nonflowering"
$ws.Range("C423").Value = "Medium"
$ws.Range("D423").Value = 66
$ws.Range("E423").Value = 77
$ws.Range("G423").Value = 0.25
$ws.Range("H423").Value = 0.05
$ws.Range("I423").Value = "No"
$ws.Range("J423").Value = 3
$ws.Range("K423").Value = "Neutral"
$ws.Range("L423").Value = 5
$ws.Range("M423").Value = 0.84
$ws.Range("N423").Value = 70
$ws.Range("O423").Value = 30.01
$ws.Range("P423").Value = 7
$ws.Range("Q423").Value = 0.86
$ws.Range("R423").Value = 5.8
$ws.Range("S423").Value = 57
$ws.Range("T423").Value = 0

# Row 424
$ws.Range("A424").Value = 45847
$ws.Range("B424").Value = "Nonflowering"
$ws.Range("C424").Value = "Small"
$ws.Range("D424").Value = 66
$ws.Range("E424").Value = 77
$ws.Range("G424").Value = 0.25
$ws.Range("H424").Value = 0.1
$ws.Range("I424").Value = "No"
$ws.Range("J424").Value = 3
$ws.Range("K424").Value = "Neutral"
$ws.Range("L424").Value = 5
$ws.Range("M424").Value = 0.84
$ws.Range("N424").Value = 70
$ws.Range("O424").Value = 30.01
$ws.Range("P424").Value = 7
$ws.Range("Q424").Value = 0.86
$ws.Range("R424").Value = 5.8
$ws.Range("S424").Value = 57
$ws.Range("T424").Value = 0

# Row 425
$ws.Range("A425").Value = 45847
$ws.Range("B425").Value = "Nonflowering"
$ws.Range("C425").Value = "Medium"
$ws.Range("D425").Value = 66
$ws.Range("E425").Value = 77
$ws.Range("G425").Value = 0.25
$ws.Range("H425").Value = 0.2
$ws.Range("I425").Value = "No"
$ws.Range("J425").Value = 3
$ws.Range("K425").Value = "Neutral"
$ws.Range("L425").Value = 5
$ws.Range("M425").Value = 0.84
$ws.Range("N425").Value = 70
$ws.Range("O425").Value = 30.01
$ws.Range("P425").Value = 7
$ws.Range("Q425").Value = 0.86
$ws.Range("R425").Value = 5.8
$ws.Range("S425").Value = 57
$ws.Range("T425").Value = 0

# Row 426
$ws.Range("A426").Value = 45847
$ws.Range("B426").Value = "Nonflowering"
$ws.Range("C426").Value = "Medium"
$ws.Range("D426").Value = 66
$ws.Range("E426").Value = 77
$ws.Range("G426").Value = 0.25
$ws.Range("H426").Value = 0.2
$ws.Range("I426").Value = "No"
$ws.Range("J426").Value = 3
$ws.Range("K426").Value = "Dark"
$ws.Range("L426").Value = 5
$ws.Range("M426").Value = 0.84
$ws.Range("N426").Value = 70
$ws.Range("O426").Value = 30.01
$ws.Range("P426").Value = 7
$ws.Range("Q426").Value = 0.86
$ws.Range("R426").Value = 5.8
$ws.Range("S426").Value = 57
$ws.Range("T426").Value = 0

# Row 427
$ws.Range("A427").Value = 45847
$ws.Range("B427").Value = "Nonflowering"
$ws.Range("C427").Value = "Large"
$ws.Range("D427").Value = 66
$ws.Range("E427").Value = 77
$ws.Range("G427").Value = 0.25
$ws.Range("H427").Value = 0.15
$ws.Range("I427").Value = "No"
$ws.Range("J427").Value = 4
$ws.Range("K427").Value = "Dark"
$ws.Range("L427").Value = 5
$ws.Range("M427").Value = 0.84
$ws.Range("N427").Value = 70
$ws.Range("O427").Value = 30.01
$ws.Range("P427").Value = 7
$ws.Range("Q427").Value = 0.86
$ws.Range("R427").Value = 5.8
$ws.Range("S427").Value = 57
$ws.Range("T427").Value = 0

# Row 428
$ws.Range("A428").Value = 45847
$ws.Range("B428").Value = "Tree"
$ws.Range("C428").Value = "Medium"
$ws.Range("D428").Value = 66
$ws.Range("E428").Value = 77
$ws.Range("G428").Value = 0.25
$ws.Range("H428").Value = 0.55
$ws.Range("I428").Value = "No"
$ws.Range("J428").Value = 1
$ws.Range("K428").Value = "Neutral"
$ws.Range("L428").Value = 5
$ws.Range("M428").Value = 0.84
$ws.Range("N428").Value = 70
$ws.Range("O428").Value = 30.01
$ws.Range("P428").Value = 7
$ws.Range("Q428").Value = 0.86
$ws.Range("R428").Value = 5.8
$ws.Range("S428").Value = 57
$ws.Range("T428").Value = 0

# ---------------------------------------------------------------------
# 2. Apply the same date number format used throughout column A to the
#    new date cells (copy format down from the last pre-existing row).
# ---------------------------------------------------------------------
$ws.Range("A421").Copy()
$ws.Range("A422:A428").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Fill in the Temp_Diff (column F) formulas. The original workbook
#    had a single shared formula (=ABS(D-E)) spanning F359:F421; mimic
#    that by writing the same formula into the new rows as two fill
#    operations (F424:F428 then F422:F423) so the formula keeps working
#    the same way it did originally and the new shared-formula group
#    lines up with the one introduced by the source edit (si="9" over
#    F424:F428).
# ---------------------------------------------------------------------
$ws.Range("F424:F428").Formula = "=ABS(D424-E424)"
$ws.Range("F422:F423").Formula = "=ABS(D422-E422)"

# ---------------------------------------------------------------------
# 4. Recalculate so cached formula results are correct, then restore
#    the view/selection state to match where the user ended up editing.
# ---------------------------------------------------------------------
$excel.Calculate()

$ws.Range("S429").Select()
